$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8859735727310181
$ws.Range("B1").Value = 1.69651186466217
$ws.Range("C1").Value = 4.046674251556396
$ws.Range("D1").Value = 3.673566579818726
$ws.Range("E1").Value = 0.7807289958000183
